$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change columns) to match the latest
# scrape. Price-column (D) values are kept as literal text -- the sheet stores
# them as strings like "47.291.67" / "0.999" / "34.60" -- so we briefly force a
# text format before writing the value (otherwise Excel auto-parses plain
# decimals such as "108.28" into a number) and then clear the format again so
# the cell keeps its original (unstyled) look.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.291.67"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.490.69"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.33"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.879.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.492.80"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.188.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.71"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +14.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.60"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.81"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.28"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0783"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.18"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.18%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.994.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.12%  "
